$d = $word.ActiveDocument

# 1. Add KeepWithNext (-> <w:keepNext/>) to the seven body paragraphs
#    that should stay on the same page as the paragraph following them
#    (a lead-in sentence right before a "Based on what you ..." / intro
#    paragraph, or right after a section heading). Indices refer to the
#    document's paragraph numbering *before* the deletion in step 2.
$keepNextTargets = @{
    9  = "responsibility to investigate"
    10 = "your lender did not meet this obligation"
    18 = "indicate you were in a vulnerable position"
    24 = "more likely to be unsuitable"
    25 = "nature of your loan was unsuitable"
    30 = "may be dealing with a responsible lending issue"
    31 = "suffering from"
}
foreach ($idx in $keepNextTargets.Keys) {
    $p = $d.Paragraphs.Item($idx)
    if ($p.Range.Text -notmatch [regex]::Escape($keepNextTargets[$idx])) {
        throw "Paragraph $idx did not match expected text: $($p.Range.Text)"
    }
    $p.Format.KeepWithNext = $true
}

# 2. Remove the stray blank paragraph sitting between "Take this report
#    with you ..." and the "{% if asked_why or ... %}" paragraph.
$blank = $d.Paragraphs.Item(6)
if ($blank.Range.Text.Trim() -ne "") {
    throw "Paragraph 6 was not blank: $($blank.Range.Text)"
}
$blank.Range.Delete()

# 3. Give the Heading 1 and Heading 3 styles a KeepWithNext paragraph
#    setting too, so those headings stay with the text that follows them.
$d.Styles("Heading 1").ParagraphFormat.KeepWithNext = $true
$d.Styles("Heading 3").ParagraphFormat.KeepWithNext = $true

Write-Output "Done"
